$wb = $excel.ActiveWorkbook

# --- "FM mkdocs table" sheet: insert a new row for "Storage node file" ---
# In the source data ('Source table'), a "Storage node file" entry already
# exists at row 41 (between "1D roughness file" row 40 and
# "**Spatial data files**" row 42). The "FM mkdocs table" sheet is missing
# the corresponding row, so insert a new row 33 (pushing the old rows
# 33-47 down to 34-48) and fill it with the same formula pattern used by
# the surrounding rows, pointed at 'Source table' row 41.

$wsFM = $wb.Worksheets.Item("FM mkdocs table")
$wsFM.Rows.Item(33).Insert()

$wsFM.Range("A33").Formula = "=IF(ISBLANK('Source table'!A41),"" "",'Source table'!A41)"
$wsFM.Range("B33").Formula = "=IFERROR(VLOOKUP('Source table'!B41,'mkdocs symbols'!`$A`$1:`$C`$5,2,0),"" "")"
$wsFM.Range("C33").Formula = "=IFERROR(VLOOKUP('Source table'!C41,'mkdocs symbols'!`$A`$1:`$C`$5,2,0),"" "")"
$wsFM.Range("D33").Formula = "=IF(ISBLANK('Source table'!D41),"" "",'Source table'!D41)"
$wsFM.Range("E33").Formula = "=IF(OR(ISBLANK('Source table'!E41),ISBLANK('Source table'!F41)),"" "",""[""&'Source table'!F41&""][""&'Source table'!E41&"".""&'Source table'!F41&""]"")"
$wsFM.Range("F33").Formula = "=IF(ISBLANK('Source table'!G41),"" "",""_""&'Source table'!G41&""_"")"

# --- Update the active sheet / selection state ---
# "Source table" was the active sheet with G4 selected; the edit leaves
# "FM mkdocs table" active instead, with A51 selected (and "Source
# table"'s remembered selection moved to A41).

$wsSource = $wb.Worksheets.Item("Source table")
$wsSource.Range("A41").Select()

$wsFM.Activate()
$wsFM.Range("A51").Select()
